$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top, pushing the existing header row down to row 2
$ws.Rows("1:1").Insert()

# Set the new cell A1 to the added text
$ws.Range("A1").Value = "Jeremy Tryon"

# Reset the selection back to the default (A1) so no stale selection is saved
$ws.Range("A1").Select() | Out-Null
